# Restore the "Artificial Intelligence and Machine Learning" content of the
# Comprehensive Budget template (it currently still carries stale
# "Finance - Core Banking System Modernization" text/numbers from before the
# last restore).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Instructions & User Guide
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning Comprehensive Budget - User Guide & Instructions"
$ws.Range("A56").Value = "📋 ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING PROJECT OVERVIEW"

# ---------------------------------------------------------------------
# 2) Budget Summary
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Executive Budget Summary"

# ---------------------------------------------------------------------
# 3) Resources
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Resources Budget"

$ws.Range("A4").Value = "Data Scientists"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 16

$ws.Range("A5").Value = "ML Engineers"
$ws.Range("B5").Value = 165
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 16

$ws.Range("A6").Value = "AI Architects"
$ws.Range("B6").Value = 200

$ws.Range("A7").Value = "DevOps Engineers"
$ws.Range("B7").Value = 150
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 16

$ws.Range("A8").Value = "Project Manager"
$ws.Range("B8").Value = 140
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 18

$ws.Range("A9").Value = "Business Analysts"
$ws.Range("B9").Value = 120
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 12

$ws.Range("A10").Value = "QA Engineers"
$ws.Range("B10").Value = 110
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 10

# ---------------------------------------------------------------------
# 4) Logistics
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Logistics Budget"

$ws.Range("B4").Value = 85000
$ws.Range("B5").Value = 120000
$ws.Range("B6").Value = 35000
$ws.Range("B7").Value = 25000
$ws.Range("B8").Value = 15000

# ---------------------------------------------------------------------
# 5) Technology
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Technology Budget"

$ws.Range("A4").Value = "Cloud Infrastructure (AWS/Azure)"
$ws.Range("B4").Value = 180000

$ws.Range("A5").Value = "ML Platform Licenses"
$ws.Range("B5").Value = 120000

$ws.Range("A6").Value = "Data Storage and Processing"
$ws.Range("B6").Value = 95000

$ws.Range("A7").Value = "Development Tools"
$ws.Range("B7").Value = 45000

$ws.Range("A8").Value = "Security and Compliance Tools"
$ws.Range("B8").Value = 35000

$ws.Range("A9").Value = "Monitoring and Analytics"
$ws.Range("B9").Value = 25000

# ---------------------------------------------------------------------
# 6) Training
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Training Budget"

$ws.Range("A4").Value = "AI/ML Certification Programs"
$ws.Range("B4").Value = 45000
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "Technical Training Materials"
$ws.Range("B5").Value = 25000
$ws.Range("C5").Value = 1

$ws.Range("A6").Value = "Conference and Workshop Attendance"
$ws.Range("B6").Value = 35000
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "Internal Training Development"
$ws.Range("B7").Value = 20000
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "External Training Consultants"
$ws.Range("B8").Value = 30000
$ws.Range("C8").Value = 1

# ---------------------------------------------------------------------
# 7) Contingency
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Contingency Budget"

$ws.Range("D6").Value = "Additional requirements or scope expansion"
$ws.Range("D8").Value = "Staff turnover or skill gaps"
$ws.Range("D9").Value = "Delays or timeline extensions"

# ---------------------------------------------------------------------
# 8) Timeline
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Budget Timeline"
